$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04745578092654933
$ws.Range("H2").Value = -1.606531576008338
$ws.Range("I2").Value = 13.57457575347232

$ws.Range("G3").Value = 0.03722383007507604
$ws.Range("H3").Value = -2.953268258159161

$ws.Range("G4").Value = -0.4600504391005563
$ws.Range("H4").Value = -1.608428779299619

$ws.Range("G5").Value = -0.4553137631116126
$ws.Range("H5").Value = 4.912175712690217

$ws.Range("G6").Value = 0.2425080272062985
$ws.Range("H6").Value = 3.799303887213148

$ws.Range("G7").Value = 0.2547089450619264
$ws.Range("H7").Value = 15.47518691705995

$ws.Range("G8").Value = 0.1707807670296438
$ws.Range("H8").Value = 2.382599196512085

$ws.Range("G9").Value = 0.16294769692268
$ws.Range("H9").Value = -5.269206885978187

$ws.Range("G10").Value = -0.01452701250123707
$ws.Range("H10").Value = -206.5105632882675

$ws.Range("G11").Value = -0.01951024782623787
$ws.Range("H11").Value = -33.12763019011324

$ws.Range("G12").Value = 0.1301882571772433
$ws.Range("H12").Value = -4.779044090063612

$ws.Range("G13").Value = 0.1316668645941095
$ws.Range("H13").Value = 5.639777605684797

$ws.Range("G14").Value = 0.2575602000003243
$ws.Range("H14").Value = 4.134791850717202

$ws.Range("G15").Value = 0.2612642532234014
$ws.Range("H15").Value = 3.398505971491379

$ws.Range("G16").Value = 0.1402298406700088
$ws.Range("H16").Value = -8.626762997191749

$ws.Range("G17").Value = 0.1363232666592981
$ws.Range("H17").Value = -9.728776904377376

$ws.Range("G18").Value = -0.008031533145252916
$ws.Range("H18").Value = 50.95218854870722

$ws.Range("G19").Value = 0.003655291200755492
$ws.Range("H19").Value = 534.0871451579419

$ws.Range("G20").Value = 0.158929440395021
$ws.Range("H20").Value = 14.61629158236386

$ws.Range("G21").Value = 0.1513624081333621
$ws.Range("H21").Value = 5.777953992082335

$ws.Range("G22").Value = 0.1822023985173604
$ws.Range("H22").Value = -2.156075514412821

$ws.Range("G23").Value = 0.1728799314658284
$ws.Range("H23").Value = -3.666864076545257

$ws.Range("G24").Value = -0.09403355025180332
$ws.Range("H24").Value = 0.3852997677378087

$ws.Range("G25").Value = -0.09494194923094661
$ws.Range("H25").Value = 4.684848648628144

$ws.Range("G26").Value = 0.2323174711249414
$ws.Range("H26").Value = 0.9459565119952299

$ws.Range("G27").Value = 0.2397303704804443
$ws.Range("H27").Value = 3.077810274404344

$ws.Range("G28").Value = 0.06298725223967681
$ws.Range("H28").Value = 7.115233701918354

$ws.Range("G29").Value = 0.06969247728541238
$ws.Range("H29").Value = -1.263795669506175
